# "Generate Report for Handoff"
#
# The handoff report was regenerated: a new handoff id (GUID) and a new
# content hash were minted, and the handoff/target timestamps moved a few
# seconds later. This updates the file/row-key text that is repeated across
# the "Overview", "zh-cn" and "de-de" sheets (each sheet's row 2 echoes the
# same handoff id, and the two localized sheets additionally carry the
# per-locale target-file name + its datetime).

$wb = $excel.ActiveWorkbook

$oldGuid = "c3b6b610-a6c9-4bee-a725-e2f60c03aef3"
$newGuid = "16e36900-1fdd-40a5-aba9-c0746dfc8c6f"
$oldHash = "78a7bf5819f1ce13ec2ca6daa9984b13579b4d34"
$newHash = "0838fc47f56d35ff555030d39fab9e6d7b6501fe"

# --- Overview sheet: handoff id + "Latest Handoff Date" bump ------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("D2").Value = "2016-45-19 00:45:59"

# --- zh-cn sheet: handoff id + target xlf name/hash + its datetime ------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-19 00:45:56"

# --- de-de sheet: handoff id + target xlf name/hash + its datetime ------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-19 00:45:59"
